# Reorder the "Periodo Mora" values in column E (rows 16-19) from
# descending (1811,1810,1809,1808) to ascending (1808,1809,1810,1811),
# matching the new batch of account-statement periods added upstream.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "1808"
$ws.Range("E17").Value = "1809"
$ws.Range("E18").Value = "1810"
$ws.Range("E19").Value = "1811"
